# account/lecturer.xlsx — add a new lecturer row, duplicate-name handling,
# a new "Mã giảng viên" (lecturer code) column, and highlight fills for
# lecturers that now appear on more than one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the duplicate-name case: "Ngo Van Q" -> "Dang Van Q" --------
$ws.Range("A5").Value = "Đặng Văn Q"

# --- 2. Insert a duplicate "Lý Thị B" row (she teaches a 2nd GPS class) -
# Shifts the old row 7 (Đoàn Văn C) down; we rebuild its content below.
$ws.Rows(8).Insert()
$ws.Range("A8").Value = "Lý Thị B"
$ws.Range("B8").Value = "GPS"
$ws.Range("C8").Value = 123

# --- 3. Add a brand-new row for "Lê Văn Y" teaching GIS_DC too ---------
$ws.Range("A10").Value = "Lê Văn Y"
$ws.Range("B10").Value = "GIS_DC"
$ws.Range("C10").Value = 123

# --- 4. New column D: "Mã giảng viên" (lecturer code) -------------------
$ws.Range("D1").Value = "Mã giảng viên"
$ws.Range("D2").Value = "n01"
$ws.Range("D3").Value = "n02"
$ws.Range("D4").Value = "n03"
$ws.Range("D5").Value = "k01"
$ws.Range("D6").Value = "k02"
$ws.Range("D7").Value = "gp01"
$ws.Range("D8").Value = "gp02"
$ws.Range("D9").Value = "gi01"
$ws.Range("D10").Value = "gi02"

# --- 5. Highlight the duplicated-name rows ------------------------------
# Yellow for the "Lê Văn Y" duplicate (rows 3 and 10)
$ws.Range("A3:C3").Interior.Color = 65535
$ws.Range("A10:C10").Interior.Color = 65535

# Orange for the "Lý Thị B" duplicate (rows 7 and 8)
$ws.Range("A7:C8").Interior.Color = 49407

# --- 6. Column D width + selection cosmetics ----------------------------
$ws.Columns("D").ColumnWidth = 12
$ws.Range("F8").Select()
